# Daily attendance processing - 2026-01-21 12:04:04
# Applies:
#   1. Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" in column G wherever found.
#   2. Mark six "Not Recorded" session rows (14, 97, 117, 137, 157, 177) as Recorded,
#      restyling columns A:I to the normal "Recorded" look, filling Recorded By / Students / Status.
#   3. Recompute the dependent summary figures (Class Statistics block, per-group stats block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the order of "System" and the email address in the "Recorded By" column
# ---------------------------------------------------------------------------
$usedRows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}

# ---------------------------------------------------------------------------
# 2) Flip the newly-recorded sessions from "Not Recorded" to "Recorded"
# ---------------------------------------------------------------------------
$newlyRecorded = @(
    @{ Row = 14;  Students = "16/28" },
    @{ Row = 97;  Students = "14/31" },
    @{ Row = 117; Students = "8/28" },
    @{ Row = 137; Students = "19/29" },
    @{ Row = 157; Students = "16/33" },
    @{ Row = 177; Students = "19/30" }
)

foreach ($item in $newlyRecorded) {
    $row = $item.Row

    # Re-use the formatting of an ordinary "Recorded" row (row 2) for columns A:I
    $ws.Range("A2:I2").Copy()
    $ws.Range("A" + $row + ":I" + $row).PasteSpecial(-4122)

    $ws.Range("G" + $row).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $row).Value = $item.Students
    $ws.Range("I" + $row).Value = "Recorded"
}

# ---------------------------------------------------------------------------
# 3) Recompute the top "Class Statistics" box (K3:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 180
$ws.Range("L7").Value = 78

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "69.8%"
$ws.Range("K4").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "79.4%"
$ws.Range("K4").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Recompute the per-group "Group Statistics" table rows for B1-1..B1-6
# ---------------------------------------------------------------------------
$groupStats = @(
    @{ Row = 15; O = 16; P = 4; R = "80.0%"; S = "79.9%" },
    @{ Row = 19; O = 16; P = 4; R = "80.0%"; S = "83.5%" },
    @{ Row = 20; O = 16; P = 4; R = "80.0%"; S = "85.5%" },
    @{ Row = 21; O = 16; P = 4; R = "80.0%"; S = "84.5%" },
    @{ Row = 22; O = 16; P = 4; R = "80.0%"; S = "86.7%" },
    @{ Row = 23; O = 16; P = 4; R = "80.0%"; S = "82.3%" }
)

foreach ($g in $groupStats) {
    $row = $g.Row

    $ws.Range("O" + $row).Value = $g.O
    $ws.Range("P" + $row).Value = $g.P

    $ws.Range("R" + $row).NumberFormat = "@"
    $ws.Range("R" + $row).Value = $g.R
    $ws.Range("K4").Copy()
    $ws.Range("R" + $row).PasteSpecial(-4122)

    $ws.Range("S" + $row).NumberFormat = "@"
    $ws.Range("S" + $row).Value = $g.S
    $ws.Range("K4").Copy()
    $ws.Range("S" + $row).PasteSpecial(-4122)
}

$wb.Save()
